$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 210.73914
$ws.Range("I5").Value = 199.85715
$ws.Range("K5").Value = 199.85715
$ws.Range("M5").Value = -84.85714999999999
$ws.Range("H62").Value = 21404
$ws.Range("I62").Value = 22297.455
$ws.Range("K62").Value = 22297.455
$ws.Range("M62").Value = -21673.455
$ws.Range("H65").Value = 21404
$ws.Range("I65").Value = 22297.455
$ws.Range("K65").Value = 111487.275
$ws.Range("M65").Value = -108367.275
$ws.Range("H135").Value = 1793.2632
$ws.Range("I135").Value = 1886.6471
$ws.Range("K135").Value = 16979.8239
$ws.Range("M135").Value = -14444.8239
$ws.Range("H138").Value = 2685.7273
$ws.Range("I138").Value = 2048.0476
$ws.Range("J138").Value = 3801.6667
$ws.Range("K138").Value = 6144.1428
$ws.Range("L138").Value = 11405.0001
$ws.Range("M138").Value = -1004.1428
$ws.Range("N138").Value = -21685.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2005.4255
$ws.Range("I2").Value = 1870.878
$ws.Range("K2").Value = 1870.878
$ws.Range("M2").Value = -1757.878
$ws.Range("H5").Value = 261
$ws.Range("I5").Value = 148.125
$ws.Range("K5").Value = 148.125
$ws.Range("M5").Value = -36.125
$ws.Range("H57").Value = 8348332
$ws.Range("I57").Value = 8348332
$ws.Range("K57").Value = 8348332
$ws.Range("M57").Value = -8347848
$ws.Range("H61").Value = 1863214
$ws.Range("I61").Value = 1863214
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1863214
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1863002
$ws.Range("H74").Value = 1986644.9
$ws.Range("I74").Value = 3050883.8
$ws.Range("J74").Value = 3290.682
$ws.Range("K74").Value = 3050883.8
$ws.Range("L74").Value = 3290.682
$ws.Range("M74").Value = -3050009.8
$ws.Range("N74").Value = -5038.682
$ws.Range("H77").Value = 1986644.9
$ws.Range("I77").Value = 3050883.8
$ws.Range("J77").Value = 3290.682
$ws.Range("K77").Value = 15254419
$ws.Range("L77").Value = 16453.41
$ws.Range("M77").Value = -15250051
$ws.Range("N77").Value = -25189.41
$ws.Range("H116").Value = 2005.4255
$ws.Range("I116").Value = 1870.878
$ws.Range("K116").Value = 1870.878
$ws.Range("M116").Value = 423.1220000000001
$ws.Range("H122").Value = 2808.6
$ws.Range("I122").Value = 1472.9166
$ws.Range("J122").Value = 4812.125
$ws.Range("K122").Value = 4418.7498
$ws.Range("L122").Value = 14436.375
$ws.Range("M122").Value = -1968.7498
$ws.Range("N122").Value = -19336.375
$ws.Range("H132").Value = 487305.12
$ws.Range("I132").Value = 530542.3
$ws.Range("J132").Value = 11695.75
$ws.Range("K132").Value = 1591626.9
$ws.Range("L132").Value = 35087.25
$ws.Range("M132").Value = -1589096.9
$ws.Range("N132").Value = -40147.25
$ws.Range("H136").Value = 1863214
$ws.Range("I136").Value = 1863214
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5589642
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5587092
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("N136").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2005.4255
$ws.Range("I3").Value = 1870.878
$ws.Range("K3").Value = 1870.878
$ws.Range("M3").Value = -1756.878
$ws.Range("H4").Value = 261
$ws.Range("I4").Value = 148.125
$ws.Range("K4").Value = 148.125
$ws.Range("M4").Value = -33.125
$ws.Range("H16").Value = 24990.5
$ws.Range("J16").Value = 24990.5
$ws.Range("L16").Value = 24990.5
$ws.Range("N16").Value = -25330.5
$ws.Range("H36").Value = 3720
$ws.Range("I36").Value = 1500
$ws.Range("J36").Value = 7050
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 7050
$ws.Range("M36").Value = -966
$ws.Range("N36").Value = -8118
$ws.Range("H86").Value = 3501.6667
$ws.Range("I86").Value = 3501.6667
$ws.Range("K86").Value = 3501.6667
$ws.Range("M86").Value = -2378.6667
$ws.Range("H89").Value = 3501.6667
$ws.Range("I89").Value = 3501.6667
$ws.Range("K89").Value = 17508.3335
$ws.Range("M89").Value = -11892.3335
$ws.Range("H94").Value = 3119.5
$ws.Range("I94").Value = 1743.6
$ws.Range("K94").Value = 1743.6
$ws.Range("M94").Value = -1292.6
$ws.Range("H107").Value = 3065.1667
$ws.Range("I107").Value = 3065.1667
$ws.Range("K107").Value = 3065.1667
$ws.Range("M107").Value = -1145.1667
$ws.Range("H134").Value = 803491.8
$ws.Range("I134").Value = 955252
$ws.Range("J134").Value = 14339
$ws.Range("K134").Value = 2865756
$ws.Range("L134").Value = 43017
$ws.Range("M134").Value = -2863221
$ws.Range("N134").Value = -48087

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1288.375
$ws.Range("I5").Value = 467.83334
$ws.Range("K5").Value = 467.83334
$ws.Range("M5").Value = -355.83334
$ws.Range("H22").Value = 1240.4166
$ws.Range("I22").Value = 998.6667
$ws.Range("J22").Value = 1482.1666
$ws.Range("K22").Value = 998.6667
$ws.Range("L22").Value = 1482.1666
$ws.Range("M22").Value = -648.6667
$ws.Range("N22").Value = -2182.1666
$ws.Range("H25").Value = 7899.25
$ws.Range("I25").Value = 7899.25
$ws.Range("K25").Value = 7899.25
$ws.Range("M25").Value = -7725.25
$ws.Range("H31").Value = 7068.467
$ws.Range("I31").Value = 1758.6666
$ws.Range("J31").Value = 12378.267
$ws.Range("K31").Value = 1758.6666
$ws.Range("L31").Value = 12378.267
$ws.Range("M31").Value = -1463.6666
$ws.Range("N31").Value = -12968.267
$ws.Range("H34").Value = 7068.467
$ws.Range("I34").Value = 1758.6666
$ws.Range("J34").Value = 12378.267
$ws.Range("K34").Value = 1758.6666
$ws.Range("L34").Value = 12378.267
$ws.Range("M34").Value = -1556.6666
$ws.Range("N34").Value = -12782.267
$ws.Range("H35").Value = 21272.637
$ws.Range("J35").Value = 23333.334
$ws.Range("L35").Value = 23333.334
$ws.Range("N35").Value = -23921.334
$ws.Range("H134").Value = 1253.3055
$ws.Range("I134").Value = 1209.5588
$ws.Range("K134").Value = 3628.6764
$ws.Range("M134").Value = -1093.6764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3258.6667
$ws.Range("I5").Value = 683.1539
$ws.Range("K5").Value = 2049.4617
$ws.Range("M5").Value = -1937.4617
$ws.Range("H7").Value = 2733567.8
$ws.Range("I7").Value = 1818460.6
$ws.Range("K7").Value = 5455381.800000001
$ws.Range("M7").Value = -5455269.800000001
$ws.Range("H80").Value = 5221
$ws.Range("I80").Value = 4963
$ws.Range("K80").Value = 14889
$ws.Range("M80").Value = -13953
$ws.Range("H83").Value = 5221
$ws.Range("I83").Value = 4963
$ws.Range("K83").Value = 44667
$ws.Range("M83").Value = -39987
$ws.Range("H129").Value = 1239.8
$ws.Range("I129").Value = 981.9231
$ws.Range("J129").Value = 1718.7142
$ws.Range("K129").Value = 2945.7693
$ws.Range("L129").Value = 5156.142599999999
$ws.Range("M129").Value = 2054.2307
$ws.Range("N129").Value = -15156.1426
$ws.Range("H132").Value = 3020.8333
$ws.Range("I132").Value = 1931.6666
$ws.Range("J132").Value = 3565.4167
$ws.Range("K132").Value = 17384.9994
$ws.Range("L132").Value = 32088.7503
$ws.Range("M132").Value = -14854.9994
$ws.Range("N132").Value = -37148.7503
$ws.Range("H135").Value = 3258.6667
$ws.Range("I135").Value = 683.1539
$ws.Range("K135").Value = 6148.3851
$ws.Range("M135").Value = -3613.3851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 8176.6
$ws.Range("J9").Value = 8826.846
$ws.Range("L9").Value = 8826.846
$ws.Range("N9").Value = -9166.846
$ws.Range("H11").Value = 303883.6
$ws.Range("J11").Value = 717414.3
$ws.Range("L11").Value = 717414.3
$ws.Range("N11").Value = -717692.3
$ws.Range("H102").Value = 2186.3572
$ws.Range("I102").Value = 1526.6316
$ws.Range("J102").Value = 3579.111
$ws.Range("K102").Value = 1526.6316
$ws.Range("L102").Value = 3579.111
$ws.Range("M102").Value = 95.36840000000007
$ws.Range("N102").Value = -6823.111
$ws.Range("H132").Value = 549008.7
$ws.Range("I132").Value = 574914.1
$ws.Range("J132").Value = 4994
$ws.Range("K132").Value = 1724742.3
$ws.Range("L132").Value = 14982
$ws.Range("M132").Value = -1722212.3
$ws.Range("N132").Value = -20042

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7608.636
$ws.Range("I40").Value = 6632.8887
$ws.Range("J40").Value = 11999.5
$ws.Range("K40").Value = 6632.8887
$ws.Range("L40").Value = 11999.5
$ws.Range("M40").Value = -6496.8887
$ws.Range("N40").Value = -12271.5
$ws.Range("H132").Value = 1050894.2
$ws.Range("I132").Value = 1332404.2
$ws.Range("K132").Value = 3997212.6
$ws.Range("M132").Value = -3994682.6
$ws.Range("H136").Value = 4490.276
$ws.Range("I136").Value = 4222.0713
$ws.Range("K136").Value = 12666.2139
$ws.Range("M136").Value = -10116.2139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5534.5386
$ws.Range("I113").Value = 3190.4
$ws.Range("J113").Value = 6999.625
$ws.Range("K113").Value = 9571.200000000001
$ws.Range("L113").Value = 20998.875
$ws.Range("M113").Value = -7401.200000000001
$ws.Range("N113").Value = -25338.875
$ws.Range("H136").Value = 10303698
$ws.Range("I136").Value = 12706661
$ws.Range("K136").Value = 38119983
$ws.Range("M136").Value = -38117433
